$d = $word.ActiveDocument

# Locate the end of the paragraph "... trong mỗi kloc: 40%" so the new
# bullet about unit tests estimation can be inserted right after it.
$rng = $d.Content
$found = $rng.Find.Execute("i kloc: 40%", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find anchor text 'i kloc: 40%' in the document."
}

$rng.Collapse(0)
$rng.InsertParagraphAfter() | Out-Null
$rng.Move(1, 1) | Out-Null

$rng.InsertAfter("- Số unit tests: 40") | Out-Null
